$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.891.24"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.703.44"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.55"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4044"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.005"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.62"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.469"
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08811"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.87"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.553"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.043"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001349"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "1.726.31"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.52"
$ws.Range("E18").Value = "  -3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07159"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.07"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.241"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.58"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "24.901.14"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.326"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("B26").Value = "HuobiToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.797"
$ws.Range("E26").Value = "  +30.35%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.898"
$ws.Range("E27").Value = "  -5.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.04"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.22"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "145.40"
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.232"
$ws.Range("E31").Value = "  -6.27%  "
$ws.Range("D32").Value = "1.928.05"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.263"
$ws.Range("E33").Value = "  +13.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08815"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03193"
$ws.Range("E35").Value = "  +8.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.339"
$ws.Range("E36").Value = "  -5.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.016"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2842"
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8434"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.79"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09377"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.11"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.91"
$ws.Range("E43").Value = "  +7.76%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.468"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.721"
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7436"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.391"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.93"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08348"
$ws.Range("E51").Value = "  +3.71%  "
